$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ucn2"
$ws.Cells.Item(2, 3).Value = "Crhr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.3617913333333333
$ws.Cells.Item(2, 8).Value = 1.085374
$ws.Cells.Item(2, 9).Value = 0.04578080100647625
$ws.Cells.Item(2, 10).Value = 0.06509043373192533
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.471861
$ws.Cells.Item(2, 14).Value = 1.415583
$ws.Cells.Item(2, 15).Value = 0.1161318617815716
$ws.Cells.Item(2, 16).Value = 0.1646379375675581
$ws.Cells.Item(2, 17).Value = 0.170715220338
$ws.Cells.Item(2, 18).Value = 1.536436983042
$ws.Cells.Item(2, 19).Value = 0.005316609654733736
$ws.Cells.Item(2, 20).Value = 0.010716354765002

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ucn2"
$ws.Cells.Item(3, 3).Value = "Crhr2"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.3617913333333333
$ws.Cells.Item(3, 8).Value = 1.085374
$ws.Cells.Item(3, 9).Value = 0.04578080100647625
$ws.Cells.Item(3, 10).Value = 0.06509043373192533
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.5912875
$ws.Cells.Item(3, 14).Value = 7.182575
$ws.Cells.Item(3, 15).Value = 0.8838681382184285
$ws.Cells.Item(3, 16).Value = 0.8353620624324419
$ws.Cells.Item(3, 17).Value = 1.299296693008333
$ws.Cells.Item(3, 18).Value = 7.795780158049999
$ws.Cells.Item(3, 19).Value = 0.04046419135174252
$ws.Cells.Item(3, 20).Value = 0.05437407896692333

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Ucn2"
$ws.Cells.Item(4, 3).Value = "Crhr2"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.2615406666666666
$ws.Cells.Item(4, 8).Value = 0.7846219999999999
$ws.Cells.Item(4, 9).Value = 0.03309515765745578
$ws.Cells.Item(4, 10).Value = 0.04705418251737255
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.471861
$ws.Cells.Item(4, 14).Value = 1.415583
$ws.Cells.Item(4, 15).Value = 0.1161318617815716
$ws.Cells.Item(4, 16).Value = 0.1646379375675581
$ws.Cells.Item(4, 17).Value = 0.123410840514
$ws.Cells.Item(4, 18).Value = 1.110697564626
$ws.Cells.Item(4, 19).Value = 0.003843402274714977
$ws.Cells.Item(4, 20).Value = 0.007746903563587666

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ucn2"
$ws.Cells.Item(5, 3).Value = "Crhr2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.2615406666666666
$ws.Cells.Item(5, 8).Value = 0.7846219999999999
$ws.Cells.Item(5, 9).Value = 0.03309515765745578
$ws.Cells.Item(5, 10).Value = 0.04705418251737255
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.5912875
$ws.Cells.Item(5, 14).Value = 7.182575
$ws.Cells.Item(5, 15).Value = 0.8838681382184285
$ws.Cells.Item(5, 16).Value = 0.8353620624324419
$ws.Cells.Item(5, 17).Value = 0.9392677269416666
$ws.Cells.Item(5, 18).Value = 5.63560636165
$ws.Cells.Item(5, 19).Value = 0.02925175538274081
$ws.Cells.Item(5, 20).Value = 0.03930727895378489

# Row 6
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Ucn2"
$ws.Cells.Item(6, 3).Value = "Crhr2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.091943
$ws.Cells.Item(6, 8).Value = 0.275829
$ws.Cells.Item(6, 9).Value = 0.01163439750796992
$ws.Cells.Item(6, 10).Value = 0.01654160616141831
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.471861
$ws.Cells.Item(6, 14).Value = 1.415583
$ws.Cells.Item(6, 15).Value = 0.1161318617815716
$ws.Cells.Item(6, 16).Value = 0.1646379375675581
$ws.Cells.Item(6, 17).Value = 0.043384315923
$ws.Cells.Item(6, 18).Value = 0.390458843307
$ws.Cells.Item(6, 19).Value = 0.001351124243307424
$ws.Cells.Item(6, 20).Value = 0.002723375922470721

# Row 7
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Ucn2"
$ws.Cells.Item(7, 3).Value = "Crhr2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.091943
$ws.Cells.Item(7, 8).Value = 0.275829
$ws.Cells.Item(7, 9).Value = 0.01163439750796992
$ws.Cells.Item(7, 10).Value = 0.01654160616141831
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.5912875
$ws.Cells.Item(7, 14).Value = 7.182575
$ws.Cells.Item(7, 15).Value = 0.8838681382184285
$ws.Cells.Item(7, 16).Value = 0.8353620624324419
$ws.Cells.Item(7, 17).Value = 0.3301937466125
$ws.Cells.Item(7, 18).Value = 1.981162479675
$ws.Cells.Item(7, 19).Value = 0.01028327326466249
$ws.Cells.Item(7, 20).Value = 0.01381823023894758

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Ucn2"
$ws.Cells.Item(8, 3).Value = "Crhr2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.1542133333333333
$ws.Cells.Item(8, 8).Value = 0.46264
$ws.Cells.Item(8, 9).Value = 0.01951403827402921
$ws.Cells.Item(8, 10).Value = 0.02774475734791688
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.471861
$ws.Cells.Item(8, 14).Value = 1.415583
$ws.Cells.Item(8, 15).Value = 0.1161318617815716
$ws.Cells.Item(8, 16).Value = 0.1646379375675581
$ws.Cells.Item(8, 17).Value = 0.07276725768
$ws.Cells.Item(8, 18).Value = 0.65490531912
$ws.Cells.Item(8, 19).Value = 0.002266201595639858
$ws.Cells.Item(8, 20).Value = 0.004567839628073388

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Ucn2"
$ws.Cells.Item(9, 3).Value = "Crhr2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.1542133333333333
$ws.Cells.Item(9, 8).Value = 0.46264
$ws.Cells.Item(9, 9).Value = 0.01951403827402921
$ws.Cells.Item(9, 10).Value = 0.02774475734791688
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.5912875
$ws.Cells.Item(9, 14).Value = 7.182575
$ws.Cells.Item(9, 15).Value = 0.8838681382184285
$ws.Cells.Item(9, 16).Value = 0.8353620624324419
$ws.Cells.Item(9, 17).Value = 0.5538244163333333
$ws.Cells.Item(9, 18).Value = 3.322946498
$ws.Cells.Item(9, 19).Value = 0.01724783667838935
$ws.Cells.Item(9, 20).Value = 0.02317691771984349

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ucn2"
$ws.Cells.Item(10, 3).Value = "Crhr2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.0331985
$ws.Cells.Item(10, 8).Value = 14.066397
$ws.Cells.Item(10, 9).Value = 0.8899756055540688
$ws.Cells.Item(10, 10).Value = 0.8435690202413669
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.471861
$ws.Cells.Item(10, 14).Value = 1.415583
$ws.Cells.Item(10, 15).Value = 0.1161318617815716
$ws.Cells.Item(10, 16).Value = 0.1646379375675581
$ws.Cells.Item(10, 17).Value = 3.3186920774085
$ws.Cells.Item(10, 18).Value = 19.912152464451
$ws.Cells.Item(10, 19).Value = 0.1033545240131756
$ws.Cells.Item(10, 20).Value = 0.1388834636884243

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Ucn2"
$ws.Cells.Item(11, 3).Value = "Crhr2"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 7.0331985
$ws.Cells.Item(11, 8).Value = 14.066397
$ws.Cells.Item(11, 9).Value = 0.8899756055540688
$ws.Cells.Item(11, 10).Value = 0.8435690202413669
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 3.5912875
$ws.Cells.Item(11, 14).Value = 7.182575
$ws.Cells.Item(11, 15).Value = 0.8838681382184285
$ws.Cells.Item(11, 16).Value = 0.8353620624324419
$ws.Cells.Item(11, 17).Value = 25.25823785806875
$ws.Cells.Item(11, 18).Value = 101.032951432275
$ws.Cells.Item(11, 19).Value = 0.7866210815408933
$ws.Cells.Item(11, 20).Value = 0.7046855565529426
